# Bond screener date-rollover update.
# "Today" reference used throughout the sheet moves from 2023-11-20 (serial 45250)
# to 2023-11-21 (serial 45251). For every data row (2..262):
#   - Column G ("Dni od poprzedniej wyplaty" / days since previous payment) = TODAY - F, so it increases by 1.
#   - Column I ("Dni do nastepnej wyplaty" / days until next payment) = H - TODAY, so it decreases by 1.
# A handful of rows (129, 201, 202, 203, 204, 251) additionally roll over a coupon
# period boundary, so F and/or H (and consequently G/I) change to new explicit values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 262

# Rows whose period boundary rolled over - handled explicitly below, skipped in the generic loop.
$specialRows = @(129, 201, 202, 203, 204, 251)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($specialRows -contains $r) {
        continue
    }

    $gCell = $ws.Cells.Item($r, 7)   # column G
    $iCell = $ws.Cells.Item($r, 9)   # column I

    if ($gCell.Value2 -ne $null) {
        $gCell.Value = $gCell.Value2 + 1
    }
    if ($iCell.Value2 -ne $null) {
        $iCell.Value = $iCell.Value2 - 1
    }
}

# Rows 129 and 251: previous payment date rolled forward to 2023-11-20 (serial 45250);
# next payment date (H) stays the same; days recompute accordingly.
foreach ($r in @(129, 251)) {
    $ws.Cells.Item($r, 6).Value = 45250   # column F
    $ws.Cells.Item($r, 7).Value = 1       # column G
    $ws.Cells.Item($r, 9).Value = 91      # column I
}

# Rows 201-204: next payment date rolled forward to 2024-05-21 (serial 45433);
# previous payment date (F) stays the same; days recompute accordingly.
foreach ($r in @(201, 202, 203, 204)) {
    $ws.Cells.Item($r, 7).Value = 184     # column G
    $ws.Cells.Item($r, 8).Value = 45433   # column H
    $ws.Cells.Item($r, 9).Value = 182     # column I
}
